$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per the crypto price/volume refresh diff.
# Values that look like plain numbers (e.g. "605.70") are prefixed with a
# leading apostrophe so Excel stores them as text, matching the original
# inlineStr/text representation of the price column instead of converting
# them to numeric cells.
$ws.Range("D2").Value = "67.746.26"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "3.820.90"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'605.70"
$ws.Range("E5").Value = "  +1.90%  "
$ws.Range("D6").Value = "'166.30"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("E10").Value = "  +1.20%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "'35.95"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "4.467.05"
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").Value = "3.852.75"
$ws.Range("E15").Value = "  +3.65%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'18.55"
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "67.773.83"
$ws.Range("E17").Value = "  +0.37%  "
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("D20").Value = "'462.29"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("D21").Value = "'9.92"
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").Value = "'0.702"
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("D23").Value = "'0.0000148"
$ws.Range("E23").Value = "  -2.43%  "
$ws.Range("D24").Value = "'83.31"
$ws.Range("D25").Value = "'12.15"
$ws.Range("E25").Value = "  +2.49%  "
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").Value = "'10.05"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("D29").Value = "3.971.49"
$ws.Range("E29").Value = "  +1.50%  "
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("D31").Value = "'7.36"
$ws.Range("E31").Value = "  +1.64%  "
$ws.Range("D32").Value = "'2.24"
$ws.Range("E32").Value = "  +1.03%  "
$ws.Range("D33").Value = "'29.62"
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "'9.12"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.100"
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "'3.30"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.138"
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").Value = "'5.82"
$ws.Range("E39").Value = "  +1.56%  "
$ws.Range("D40").Value = "'0.996"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D43").Value = "'44.89"
$ws.Range("E43").Value = "  -2.06%  "
$ws.Range("D44").Value = "'28.91"
$ws.Range("E44").Value = "  +10.15%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'47.66"
$ws.Range("E45").Value = "  -2.02%  "
$ws.Range("D46").Value = "'1.43"
$ws.Range("E46").Value = "  +16.03%  "
$ws.Range("B47").Value = "TheGraph"
$ws.Range("C47").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D47").Value = "'0.301"
$ws.Range("E47").Value = "  +0.84%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'150.93"
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "'8.36"
$ws.Range("E49").Value = "  +0.73%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'1.86"
$ws.Range("E50").Value = "  +2.34%  "
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "'387.45"
$ws.Range("E51").Value = "  -0.57%  "
